# Update the "Förändrad" (Changed) date column (C) for rows 2-14
# from 2023-10-09 (serial 45208) to 2023-10-13 (serial 45212).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 3).Value = 45212
}
